$wb = $excel.ActiveWorkbook

# The Overview rollup mirrors the per-language handoff status text, so it
# needs to flip in lockstep with the zh-cn / de-de sheets below.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Handoff failed"
$overview.Range("C2").Value = "Handoff failed"

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Remove the hyperlink on C2 (the handoff-file link) without touching
    # the other hyperlinks on the sheet (A2, A3).
    $linksToRemove = @()
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq '$C$2') {
            $linksToRemove += $hl
        }
    }
    foreach ($hl in $linksToRemove) {
        $hl.Delete()
    }

    # Fully clear C2 (value + formatting) now that its hyperlink is gone.
    $ws.Range("C2").Clear()

    # Report generated for a failed handoff instead of "not yet handed off".
    $ws.Range("B2").Value = "Handoff failed"

    # The handoff never produced a target file / date, so the "latest
    # handoff datetime" collapses to the zero-date placeholder.
    $ws.Range("D2").Value = "0001-01-01 00:00:00"

    # Handoff reason is now "Ignored" rather than "Include".
    $ws.Range("H2").Value = "Ignored"
}
